# Append "+526" to every formula in column A (rows 1-100). The original
# formulas are B{r}*C{r}+SUM(D{r},E{r},F{r})/G{r}-H{r}; the new formulas
# are the same expression plus 526, which changes every computed result
# from 2 to 528. Rows 2-65 and 66-100 are stored as Excel "shared
# formula" groups in the saved XML; rewriting every cell with the same
# formula pattern lets that grouping be reconstructed on save exactly as
# before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 100; $r++) {
    $ws.Cells.Item($r, 1).Formula = "=B$r*C$r+SUM(D$r,E$r,F$r)/G$r-H$r+526"
}

$excel.Calculate()
